$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 summary stats ("HIGH","LOW","WIDTH","AVE","STD" over the pe column) ---
# HIGH/LOW/WIDTH (B4:D4) are unaffected by the new data; only AVE (E4) and STD (F4) move.
$ws.Range("E4").Value = 120.506067891204
$ws.Range("F4").Value = 19.56020479652813

# --- Rows 6-248: ave_pe / plus_1sigma / minus_1sigma (columns D:F) recalculated for the
#     whole historical series now that two new trading days were added ---
$firstRow = 6
$lastRow = 248
$n = $lastRow - $firstRow + 1
$defArr = New-Object 'object[,]' $n,3
for ($i = 0; $i -lt $n; $i++) {
    $defArr[$i,0] = 120.506067891204
    $defArr[$i,1] = 140.0662726877321
    $defArr[$i,2] = 100.9458630946758
}
$ws.Range("D6:F248").Value = $defArr

# --- Row 249: was an accidental duplicate of row 247 (date 44434) - replace with the
#     correct next trading day (44438) and its real pe / per-stock figures ---
$row249 = @(44438,135.4510970710862,120.506067891204,140.0662726877321,100.9458630946758,1.699885591071942,1.265303684032546,0.4278285423837578,41.07086885957026,-0.0268226945223286,0.04610367011235595,2.187190312172779,0.2164180581641915,0.3491323825976512,0.5227842860309949,1.579640214992065,8.555855532054183,0.05527887733600029,0.9017270744579109,3.317535251154265,3.196678492209704,0.2492445382094148,4.148578161504502,42.02890257998334,0.2510848981518805,3.06216968826586,1.302318124229708,1.416805917576391,0.6973113663913429,0.25339575551732,1.497517587356535,3.404451058206868,0.06498674022675001,0.06534368476545728,7.595230225107561,1.258104476714067,0.4869123229152642,0.4276456253391171,0.1052511817817716,0.3231538470983296,-0.1036793249187548,0.6162597420281106,0.9347007408170956)
$arr249 = New-Object 'object[,]' 1,$row249.Length
for ($i = 0; $i -lt $row249.Length; $i++) { $arr249[0,$i] = $row249[$i] }
$ws.Range("B249:AR249").Value = $arr249

# --- Row 250: was an accidental duplicate of row 248 (date 44435) - replace with the
#     correct next trading day (44439) and its real pe / per-stock figures ---
$row250 = @(44439,133.0531027868114,120.506067891204,140.0662726877321,100.9458630946758,1.656092144044958,1.371819868266602,0.4393872256570433,40.2278934218232,-0.02709805253473216,0.04715424912121458,2.230934118379613,0.2001730064465059,0.3404067550974208,0.5062660598542961,1.53029880487951,8.477802113131741,0.05216367895596714,0.8965648798000071,3.219148409485037,3.153273472799731,0.236309369170044,4.151931163779998,41.28810162330976,0.2345602517973264,2.953222314854359,1.35212179584781,1.38082028106722,0.7006883311564949,0.2429424546667041,1.53972121791823,3.307181027673144,0.0642201862718153,0.06487610919793697,7.212784886897905,1.232717812957244,0.4681401502990308,0.4380516689333572,0.1058231136044946,0.3190041989955765,-0.1047092520680748,0.5331524882149147,1.009161437058023)
$arr250 = New-Object 'object[,]' 1,$row250.Length
for ($i = 0; $i -lt $row250.Length; $i++) { $arr250[0,$i] = $row250[$i] }
$ws.Range("B250:AR250").Value = $arr250
